# aggiornamento fino a 1/09/2021
# Append 9 new daily rows (358-366, dates 2021-08-24 .. 2021-09-01) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
    @(44432, 2, 9, 259.8902685532775),
    @(44433, 0, 9, 259.8902685532775),
    @(44434, 0, 8, 231.0135720473578),
    @(44435, 1, 4, 115.5067860236789),
    @(44436, 2, 6, 173.2601790355183),
    @(44437, 1, 6, 173.2601790355183),
    @(44438, 0, 6, 173.2601790355183),
    @(44439, 1, 5, 144.3834825295986),
    @(44440, 0, 5, 144.3834825295986)
)

$startRow = 358

# Copy formatting (date style with border/alignment on column A) from the
# last existing data row down into each new row before filling in values,
# so the new rows match the look of the existing table.
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $srcRow = $ws.Range("A357:D357")
    $dstRow = $ws.Range("A" + $r + ":D" + $r)
    $srcRow.Copy($dstRow)
}

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
